{"js": "// Rename the \"estado\" field (String) to \"aceptado\" (Boolean) in the\n// \"Equipo\" collection documentation, per commit: Create new field\n// \"Aceptado\" in Equipo and Jugador collection.\n//\n// Strategy: locate the paragraph that starts the \"estado: String (...\"\n// field description (it spans three paragraphs in the original), build\n// a Range covering those three paragraphs, and replace that range's\n// content with the new two-paragraph OOXML describing the \"aceptado:\n// Boolean (...)\" field.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that contains the \"estado: \" field label (start of\n// the block to be replaced) and the paragraph that ends with the final\n// closing parenthesis of that same field description (\"a la espera de\n// serlo.)\" \u2014 end of the block to be replaced).\nlet startIndex = -1;\nlet endIndex = -1;\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (startIndex === -1 && /^\\s*estado:\\s*String\\s*\\(/.test(t)) {\n    startIndex = i;\n  }\n  if (startIndex !== -1 && i >= startIndex && /a la espera de serlo\\.\\)\\s*$/.test(t)) {\n    endIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1 || endIndex === -1) {\n  throw new Error(\"Could not locate the 'estado' field paragraphs to replace.\");\n}\n\nconst startRange = items[startIndex].getRange(\"Start\");\nconst endRange = items[endIndex].getRange(\"End\");\nconst targetRange = startRange.expandTo(endRange);\n\nconst newOoxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:tab/><w:t>aceptado</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t>Boolean</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t xml:space=\"preserve\">indica si el equipo ha sido aceptado en la liga o se </w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:u w:val=\"single\"/><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:tab/><w:t xml:space=\"preserve\">     </w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t>encuentra en espera de ello)</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntargetRange.insertOoxml(newOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Rename the \"estado\" field (String) to \"aceptado\" (Boolean) in the\n# \"Equipo\" collection documentation, per commit: Create new field\n# \"Aceptado\" in Equipo and Jugador collection.\n#\n# Strategy: locate the paragraph that starts the \"estado: String (...\"\n# field description (it spans three paragraphs in the original), build\n# a Range covering those three paragraphs, and replace that range's\n# content with the new two-paragraph OOXML describing the \"aceptado:\n# Boolean (...)\" field.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$startIndex = -1\n$endIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($startIndex -eq -1 -and $t -match \"estado:\\s*String\\s*\\(\") {\n        $startIndex = $i\n    }\n    if ($startIndex -ne -1 -and $i -ge $startIndex -and $t -match \"a la espera de serlo\\.\\)\") {\n        $endIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1 -or $endIndex -eq -1) {\n    throw \"Could not locate the 'estado' field paragraphs to replace.\"\n}\n\n$startRange = $d.Paragraphs.Item($startIndex).Range\n$endRange = $d.Paragraphs.Item($endIndex).Range\n$targetRange = $d.Range($startRange.Start, $endRange.End)\n\n$newXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:tab/><w:t>aceptado</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t>Boolean</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t xml:space=\"preserve\">indica si el equipo ha sido aceptado en la liga o se </w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:u w:val=\"single\"/><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:tab/><w:t xml:space=\"preserve\">     </w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr><w:t>encuentra en espera de ello)</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n$targetRange.InsertXML($newXml)\n"}
